# Update countries & provincias Spain
# Applies the "paises.xlsx" data refresh: updated timestamp, refreshed
# case counts for a handful of countries, and three countries (Cabo Verde,
# Timor Oriental, Fiyi) whose case totals moved them ahead of their
# previous alphabetical/ranking neighbours in the sorted table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: "last updated" timestamp -------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 19 de Abril de 2020 a las 13:52"

# --- Row 8: Alemania (Germany) refreshed counts ---------------------------
$ws.Range("B8").Value = 143779
$ws.Range("C8").Value = 55
$ws.Range("E8").Value = 51236
$ws.Range("G8").Value = 5
$ws.Range("H8").Value = 4543

# --- Row 73: Eslovenia (Slovenia) refreshed counts ------------------------
$ws.Range("B73").Value = 1330
$ws.Range("C73").Value = 13
$ws.Range("D73").Value = 192
$ws.Range("E73").Value = 1064
$ws.Range("F73").Value = 26
$ws.Range("G73").Value = 4
$ws.Range("H73").Value = 74

# --- Row 74: Lituania (Lithuania) refreshed counts ------------------------
$ws.Range("E74").Value = 1022
$ws.Range("G74").Value = 1
$ws.Range("H74").Value = 34

# --- Rows 153-155: Cabo Verde overtakes Islas Caimanes / Zambia ----------
# Cabo Verde gets fresh data and takes row 153; Islas Caimanes and Zambia
# each shift down one row, keeping their existing data.
$ws.Range("A153").Value = "Cabo Verde"
$ws.Range("C153").Value = 3
$ws.Range("D153").Value = 1
$ws.Range("E153").Value = 59
$ws.Range("F153").Value = 0

$ws.Range("A154").Value = "Islas Caimanes"
$ws.Range("C154").Value = 0
$ws.Range("D154").Value = 7
$ws.Range("E154").Value = 53
$ws.Range("F154").Value = 3
$ws.Range("H154").Value = 1

$ws.Range("A155").Value = "Zambia"
$ws.Range("B155").Value = 61
$ws.Range("C155").Value = 4
$ws.Range("D155").Value = 33
$ws.Range("E155").Value = 26
$ws.Range("F155").Value = 1
$ws.Range("H155").Value = 2

# --- Rows 180-181: Timor Oriental overtakes Laos --------------------------
$ws.Range("A180").Value = "Timor Oriental"
$ws.Range("C180").Value = 1
$ws.Range("D180").Value = 1
$ws.Range("E180").Value = 18

$ws.Range("A181").Value = "Laos"
$ws.Range("B181").Value = 19
$ws.Range("D181").Value = 2

# --- Rows 184-185: Fiyi overtakes Islas Virgenes de los Estados Unidos ---
# Case counts are tied, so only the country names swap.
$ws.Range("A184").Value = "Fiyi"
$ws.Range("A185").Value = "Islas Virgenes de los Estados Unidos"
